$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 77; this shifts the existing rows 77-100 down to 78-101.
$ws.Rows.Item(77).Insert()

# Populate the newly inserted row 77 with the new record's data.
$ws.Cells.Item(77, 1).Value = 5
$ws.Cells.Item(77, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(77, 3).Value = "Maule"
$ws.Cells.Item(77, 4).Value = 44588
$ws.Cells.Item(77, 5).Value = 7
$ws.Cells.Item(77, 6).Value = 100112030
$ws.Cells.Item(77, 7).Value = "Poroto granado"
$ws.Cells.Item(77, 8).Value = "Sin especificar"
$ws.Cells.Item(77, 9).Value = "Primera"
$ws.Cells.Item(77, 10).Value = 300
$ws.Cells.Item(77, 11).Value = 28000
$ws.Cells.Item(77, 12).Value = 28000
$ws.Cells.Item(77, 13).Value = 28000
$ws.Cells.Item(77, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(77, 15).Value = "Región del Maule"
$ws.Cells.Item(77, 16).Value = 1120
$ws.Cells.Item(77, 17).Value = 25
$ws.Cells.Item(77, 18).Value = "Hortaliza"
